# Applies the "upload samples xlsx files" commit to travel-test.xlsx
# - renames sheet1 "test" -> "01 - Profiling Questions"
# - inserts a "cover_type" column before travel_type
# - inserts "adult_input"/"child_input"/"friend_input" columns before num_adults
# - adds "plan_select"/"1-last-name"/"2-last-name" columns after get_quote_link
# - fills in row2 sample data for the new annual-trip / 2-adult scenario
# - draws a thin outline box around the adult/child/friend input block
# - narrows the url column and autofits the rest

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "01 - Profiling Questions"

# ---------------------------------------------------------------------------
# 1) Insert "cover_type" column before the old column F (travel_type)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).Insert() | Out-Null
$ws.Range("F1").Value = "cover_type"
$ws.Range("F2").Value = "ANNUAL"
$ws.Range("F1").Style = "Normal"
$ws.Range("F1:F2").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 2) Insert three columns ("adult_input", "child_input", "friend_input")
#    before what is now num_adults (column J after step 1)
# ---------------------------------------------------------------------------
$ws.Range("J1:L1").EntireColumn.Insert() | Out-Null

$ws.Range("J1:L2").NumberFormat = "@"

$ws.Range("J1").Value = "adult_input"
$ws.Range("K1").Value = "child_input"
$ws.Range("L1").Value = "friend_input"

$ws.Range("J2").Value = "Yes"
$ws.Range("K2").Value = "Yes"
$ws.Range("L2").Value = "No"

# ---------------------------------------------------------------------------
# 3) Update the (now shifted) existing columns with their new sample values
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "annual-trip"
$ws.Range("H2").Value = "01/03/2020"
$ws.Range("I2").ClearContents() | Out-Null

$ws.Range("M2").Value = "2 adults"
$ws.Range("N2").Value = "1 child"
$ws.Range("O2").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 4) Append the new trailing columns after get_quote_link (now column P)
# ---------------------------------------------------------------------------
$ws.Range("Q1").Value = "plan_select"
$ws.Range("R1").Value = "1-last-name"
$ws.Range("S1").Value = "2-last-name"

$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = "Toby"
$ws.Range("S2").Value = "Toooby"

$ws.Range("R1:S2").NumberFormat = "@"
$ws.Range("Q1").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 5) Draw a thin outline border around the adult/child/friend input block
# ---------------------------------------------------------------------------
$ws.Range("J1:O1").Borders.Item(8).LineStyle = 1
$ws.Range("J1:O1").Borders.Item(8).Weight = 2

$ws.Range("J1:J2").Borders.Item(7).LineStyle = 1
$ws.Range("J1:J2").Borders.Item(7).Weight = 2

$ws.Range("O1:O2").Borders.Item(10).LineStyle = 1
$ws.Range("O1:O2").Borders.Item(10).Weight = 2

# ---------------------------------------------------------------------------
# 6) Column widths: the other bestFit columns keep the width that shifted
#    over automatically when the new columns were inserted; only the url
#    column was deliberately narrowed (it's no longer bestFit).
# ---------------------------------------------------------------------------
$ws.Range("C1:C2").EntireColumn.ColumnWidth = 6.93

# ---------------------------------------------------------------------------
# 7) Selection / view cosmetics
# ---------------------------------------------------------------------------
$ws.Range("P7").Select() | Out-Null
